$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 3486.375
$ws.Range("I20").Value = 413
$ws.Range("J20").Value = 25000
$ws.Range("K20").Value = 413
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = -183
$ws.Range("N20").Value = -25460
$ws.Range("H35").Value = 3486.375
$ws.Range("I35").Value = 413
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 413
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -34
$ws.Range("N35").Value = -25758
$ws.Range("H137").Value = 2675
$ws.Range("I137").Value = 2375
$ws.Range("K137").Value = 7125
$ws.Range("M137").Value = -4575
$ws.Range("H138").Value = 3294.6155
$ws.Range("I138").Value = 2984.6365
$ws.Range("J138").Value = 4999.5
$ws.Range("K138").Value = 8953.9095
$ws.Range("L138").Value = 14998.5
$ws.Range("M138").Value = -3813.9095
$ws.Range("N138").Value = -25278.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10010.5
$ws.Range("J2").Value = 10010.5
$ws.Range("L2").Value = 10010.5
$ws.Range("N2").Value = -10236.5
$ws.Range("H32").Value = 3951.6667
$ws.Range("I32").Value = 2139.625
$ws.Range("K32").Value = 2139.625
$ws.Range("M32").Value = -1852.625
$ws.Range("H74").Value = 1885.4286
$ws.Range("I74").Value = 1885.4286
$ws.Range("K74").Value = 1885.4286
$ws.Range("M74").Value = -1011.4286
$ws.Range("H77").Value = 1885.4286
$ws.Range("I77").Value = 1885.4286
$ws.Range("K77").Value = 9427.143
$ws.Range("M77").Value = -5059.143
$ws.Range("H88").Value = 5201.1
$ws.Range("J88").Value = 5375.75
$ws.Range("L88").Value = 5375.75
$ws.Range("N88").Value = -6187.75
$ws.Range("H91").Value = 5201.1
$ws.Range("J91").Value = 5375.75
$ws.Range("L91").Value = 5375.75
$ws.Range("N91").Value = -8183.75
$ws.Range("H116").Value = 10010.5
$ws.Range("J116").Value = 10010.5
$ws.Range("L116").Value = 10010.5
$ws.Range("N116").Value = -14598.5
$ws.Range("H130").Value = 49995
$ws.Range("J130").Value = 49995
$ws.Range("L130").Value = 49995
$ws.Range("N130").Value = -60035
$ws.Range("H132").Value = 891.7
$ws.Range("I132").Value = 864.625
$ws.Range("K132").Value = 2593.875
$ws.Range("M132").Value = -63.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10010.5
$ws.Range("J3").Value = 10010.5
$ws.Range("L3").Value = 10010.5
$ws.Range("N3").Value = -10238.5
$ws.Range("H86").Value = 5277.8667
$ws.Range("J86").Value = 6790
$ws.Range("L86").Value = 6790
$ws.Range("N86").Value = -9036
$ws.Range("H89").Value = 5277.8667
$ws.Range("J89").Value = 6790
$ws.Range("L89").Value = 33950
$ws.Range("N89").Value = -45182
$ws.Range("H94").Value = 3325.1
$ws.Range("I94").Value = 2750.1428
$ws.Range("J94").Value = 4666.6665
$ws.Range("K94").Value = 2750.1428
$ws.Range("L94").Value = 4666.6665
$ws.Range("M94").Value = -2299.1428
$ws.Range("N94").Value = -5568.6665
$ws.Range("H105").Value = 2152.0833
$ws.Range("I105").Value = 1432.5
$ws.Range("K105").Value = 1432.5
$ws.Range("M105").Value = 314.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1099
$ws.Range("J16").Value = 1200
$ws.Range("L16").Value = 1200
$ws.Range("N16").Value = -1774
$ws.Range("H31").Value = 1896.7
$ws.Range("I31").Value = 1495
$ws.Range("K31").Value = 1495
$ws.Range("M31").Value = -1200
$ws.Range("H34").Value = 1896.7
$ws.Range("I34").Value = 1495
$ws.Range("K34").Value = 1495
$ws.Range("M34").Value = -1293
$ws.Range("H45").Value = 10172
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H58").Value = 1407.6666
$ws.Range("I58").Value = 1100.7
$ws.Range("K58").Value = 1100.7
$ws.Range("M58").Value = -897.7
$ws.Range("H113").Value = 1099
$ws.Range("J113").Value = 1200
$ws.Range("L113").Value = 1200
$ws.Range("N113").Value = -5540
$ws.Range("H132").Value = 5887.125
$ws.Range("I132").Value = 5887.125
$ws.Range("K132").Value = 17661.375
$ws.Range("M132").Value = -15131.375
$ws.Range("H136").Value = 1407.6666
$ws.Range("I136").Value = 1100.7
$ws.Range("K136").Value = 3302.1
$ws.Range("M136").Value = -752.1000000000004

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1753.5
$ws.Range("I129").Value = 1239.6
$ws.Range("J129").Value = 2039
$ws.Range("K129").Value = 3718.8
$ws.Range("L129").Value = 6117
$ws.Range("M129").Value = 1281.2
$ws.Range("N129").Value = -16117
$ws.Range("H131").Value = 1080.6923
$ws.Range("I131").Value = 1100
$ws.Range("J131").Value = 1079.0834
$ws.Range("K131").Value = 3300
$ws.Range("L131").Value = 3237.2502
$ws.Range("M131").Value = 1740
$ws.Range("N131").Value = -13317.2502
$ws.Range("H132").Value = 491.66666
$ws.Range("J132").Value = 487.5
$ws.Range("L132").Value = 4387.5
$ws.Range("N132").Value = -9447.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H47").Value = 25500
$ws.Range("J47").Value = 25500
$ws.Range("L47").Value = 25500
$ws.Range("N47").Value = -26636
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H102").Value = 975.4545000000001
$ws.Range("I102").Value = 790.8333
$ws.Range("K102").Value = 790.8333
$ws.Range("M102").Value = 831.1667
$ws.Range("H132").Value = 4804.6
$ws.Range("I132").Value = 4007.6667
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 12023.0001
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -9493.000100000001
$ws.Range("N132").Value = -23060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6077.6665
$ws.Range("I122").Value = 5794
$ws.Range("K122").Value = 17382
$ws.Range("M122").Value = -14932
$ws.Range("H132").Value = 4368.077
$ws.Range("I132").Value = 4809.778
$ws.Range("K132").Value = 14429.334
$ws.Range("M132").Value = -11899.334
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17219.25
$ws.Range("J41").Value = 17219.25
$ws.Range("L41").Value = 17219.25
$ws.Range("N41").Value = -17999.25
$ws.Range("H45").Value = 10626
$ws.Range("J45").Value = 10626
$ws.Range("L45").Value = 10626
$ws.Range("N45").Value = -11608
$ws.Range("H100").Value = 3321032.2
$ws.Range("I100").Value = 5363387.5
$ws.Range("K100").Value = 10726775
$ws.Range("M100").Value = -10726234
$ws.Range("H132").Value = 1481
$ws.Range("I132").Value = 352.5
$ws.Range("K132").Value = 1057.5
$ws.Range("M132").Value = 1472.5
$ws.Range("H137").Value = 90000
$ws.Range("J137").Value = 90000
$ws.Range("L137").Value = 90000
$ws.Range("N137").Value = -100200
